$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "Disaster"
$ws.Range("I2").Value = 0.25
$ws.Range("I3").Value = 0.05
$ws.Range("I4").Value = 0.1
$ws.Range("I5").Value = 0
$ws.Range("I6").Value = 0.15

$ws.Range("I1:I6").Font.Color = 0

$ws.Range("I7").Style = "Comma"
$ws.Range("I7").Font.ThemeColor = 1
$ws.Range("I7").Font.Size = 12
$ws.Range("I7").HorizontalAlignment = -4108

$ws.Range("I11").Select()
